$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 12 corresponds to file
# 569d02c7-ed69-4da3-bcea-4a677ba8dd86...zh-cn.xlf
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D12").Value = "2016-03-03 07:42:45"
$wsZh.Range("G12").Value = "2016-03-03 07:43:34"

# de-de sheet: row 12 corresponds to file
# 569d02c7-ed69-4da3-bcea-4a677ba8dd86...de-de.xlf
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D12").Value = "2016-03-03 07:42:56"
$wsDe.Range("G12").Value = "2016-03-03 07:43:52"
